{"js": "const replacements = [\n  [\"2023-03-29 Wednesday\", \"2023-03-30 Thursday\"],\n  [\"2+54=\", \"80-24=\"],\n  [\"52+15=\", \"41+5=\"],\n  [\"39+34=\", \"61+2=\"],\n  [\"0+54=\", \"84-38=\"],\n  [\"19+40=\", \"20+40=\"],\n  [\"6+9=\", \"41+7=\"],\n  [\"43-34=\", \"75-21=\"],\n  [\"74-51=\", \"6+32=\"],\n  [\"39+10=\", \"73-45=\"],\n  [\"20+14=\", \"75-50=\"],\n  [\"45+45=\", \"43-23=\"],\n  [\"17+61=\", \"16+63=\"],\n  [\"48-5=\", \"34+35=\"],\n  [\"76-70=\", \"87-59=\"],\n  [\"53-12=\", \"76-25=\"],\n  [\"90+8=\", \"70-35=\"],\n  [\"94-87=\", \"73-44=\"],\n  [\"33+50=\", \"70-67=\"],\n  [\"93-66=\", \"5+85=\"],\n  [\"28+17=\", \"35+44=\"],\n  [\"66+14=\", \"88-79=\"],\n  [\"5+35=\", \"47+46=\"],\n  [\"55+8=\", \"96-31=\"],\n  [\"4+70=\", \"69-48=\"],\n  [\"14+14=\", \"46-10=\"],\n  [\"12+56=\", \"88-32=\"],\n  [\"57+19=\", \"88-9=\"],\n  [\"40+28=\", \"59+25=\"],\n  [\"25-11=\", \"78-41=\"],\n  [\"42-39=\", \"47-42=\"],\n  [\"94-27=\", \"99-53=\"],\n  [\"63-61=\", \"13+50=\"],\n  [\"63-36=\", \"97-85=\"],\n  [\"19+53=\", \"92-42=\"],\n  [\"56+15=\", \"92-86=\"],\n  [\"37-28=\", \"72-65=\"],\n  [\"28+40=\", \"97-80=\"],\n  [\"69-14=\", \"94-39=\"],\n  [\"56-32=\", \"53+40=\"],\n  [\"39-10=\", \"93-53=\"],\n  [\"7+50=\", \"25+72=\"],\n  [\"82-76=\", \"10+57=\"],\n  [\"84+13=\", \"90+2=\"],\n  [\"60+12=\", \"58-0=\"],\n  [\"67-58=\", \"94-79=\"],\n  [\"10+70=\", \"94-25=\"],\n  [\"63-34=\", \"64-7=\"],\n  [\"2+11=\", \"21-7=\"],\n  [\"7+76=\", \"80+18=\"],\n  [\"76-69=\", \"3+1=\"],\n  [\"81-59=\", \"83-74=\"],\n  [\"25-0=\", \"44-12=\"],\n  [\"7+69=\", \"90+5=\"],\n  [\"35+40=\", \"88-67=\"],\n  [\"48-26=\", \"13+16=\"],\n  [\"96-68=\", \"81-79=\"],\n  [\"43+14=\", \"44+26=\"],\n  [\"43+50=\", \"99-71=\"],\n  [\"7+91=\", \"50+29=\"],\n  [\"59-28=\", \"61+0=\"],\n  [\"3+56=\", \"91-83=\"],\n  [\"1+94=\", \"51+41=\"],\n  [\"14+0=\", \"35+24=\"],\n  [\"25+4=\", \"14+8=\"],\n  [\"50-4=\", \"61+18=\"],\n  [\"4+84=\", \"7+35=\"],\n  [\"19+59=\", \"46+12=\"],\n  [\"73-17=\", \"23-15=\"],\n  [\"56-41=\", \"23+38=\"],\n  [\"18+25=\", \"78-59=\"],\n  [\"68-24=\", \"81+17=\"],\n  [\"99-98=\", \"75-10=\"],\n  [\"82-57=\", \"91-5=\"],\n  [\"73-4=\", \"14+10=\"],\n  [\"84-37=\", \"79-37=\"],\n  [\"23+35=\", \"94-22=\"],\n  [\"99-90=\", \"47+40=\"],\n  [\"69-1=\", \"17+10=\"],\n  [\"76+19=\", \"85+8=\"],\n  [\"47-43=\", \"33+51=\"],\n  [\"91-69=\", \"49+2=\"],\n  [\"73-23=\", \"30+46=\"],\n  [\"50+11=\", \"68+10=\"],\n  [\"45-9=\", \"29+66=\"],\n  [\"40+11=\", \"32+10=\"],\n  [\"8+58=\", \"50-5=\"],\n  [\"83+6=\", \"61+30=\"],\n  [\"21+48=\", \"28-1=\"],\n  [\"8+91=\", \"52-26=\"],\n  [\"17+41=\", \"36+43=\"],\n  [\"44+15=\", \"78-68=\"],\n  [\"89+1=\", \"47-7=\"],\n  [\"45-25=\", \"10+51=\"],\n  [\"98-82=\", \"24-10=\"],\n  [\"70-34=\", \"85-6=\"],\n  [\"27+52=\", \"28-10=\"],\n  [\"16+56=\", \"48+14=\"],\n  [\"49+38=\", \"88-28=\"],\n  [\"85-18=\", \"75-5=\"],\n  [\"8+82=\", \"31+28=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2023-03-29 Wednesday', '2023-03-30 Thursday'),\n    @('2+54=', '80-24='),\n    @('52+15=', '41+5='),\n    @('39+34=', '61+2='),\n    @('0+54=', '84-38='),\n    @('19+40=', '20+40='),\n    @('6+9=', '41+7='),\n    @('43-34=', '75-21='),\n    @('74-51=', '6+32='),\n    @('39+10=', '73-45='),\n    @('20+14=', '75-50='),\n    @('45+45=', '43-23='),\n    @('17+61=', '16+63='),\n    @('48-5=', '34+35='),\n    @('76-70=', '87-59='),\n    @('53-12=', '76-25='),\n    @('90+8=', '70-35='),\n    @('94-87=', '73-44='),\n    @('33+50=', '70-67='),\n    @('93-66=', '5+85='),\n    @('28+17=', '35+44='),\n    @('66+14=', '88-79='),\n    @('5+35=', '47+46='),\n    @('55+8=', '96-31='),\n    @('4+70=', '69-48='),\n    @('14+14=', '46-10='),\n    @('12+56=', '88-32='),\n    @('57+19=', '88-9='),\n    @('40+28=', '59+25='),\n    @('25-11=', '78-41='),\n    @('42-39=', '47-42='),\n    @('94-27=', '99-53='),\n    @('63-61=', '13+50='),\n    @('63-36=', '97-85='),\n    @('19+53=', '92-42='),\n    @('56+15=', '92-86='),\n    @('37-28=', '72-65='),\n    @('28+40=', '97-80='),\n    @('69-14=', '94-39='),\n    @('56-32=', '53+40='),\n    @('39-10=', '93-53='),\n    @('7+50=', '25+72='),\n    @('82-76=', '10+57='),\n    @('84+13=', '90+2='),\n    @('60+12=', '58-0='),\n    @('67-58=', '94-79='),\n    @('10+70=', '94-25='),\n    @('63-34=', '64-7='),\n    @('2+11=', '21-7='),\n    @('7+76=', '80+18='),\n    @('76-69=', '3+1='),\n    @('81-59=', '83-74='),\n    @('25-0=', '44-12='),\n    @('7+69=', '90+5='),\n    @('35+40=', '88-67='),\n    @('48-26=', '13+16='),\n    @('96-68=', '81-79='),\n    @('43+14=', '44+26='),\n    @('43+50=', '99-71='),\n    @('7+91=', '50+29='),\n    @('59-28=', '61+0='),\n    @('3+56=', '91-83='),\n    @('1+94=', '51+41='),\n    @('14+0=', '35+24='),\n    @('25+4=', '14+8='),\n    @('50-4=', '61+18='),\n    @('4+84=', '7+35='),\n    @('19+59=', '46+12='),\n    @('73-17=', '23-15='),\n    @('56-41=', '23+38='),\n    @('18+25=', '78-59='),\n    @('68-24=', '81+17='),\n    @('99-98=', '75-10='),\n    @('82-57=', '91-5='),\n    @('73-4=', '14+10='),\n    @('84-37=', '79-37='),\n    @('23+35=', '94-22='),\n    @('99-90=', '47+40='),\n    @('69-1=', '17+10='),\n    @('76+19=', '85+8='),\n    @('47-43=', '33+51='),\n    @('91-69=', '49+2='),\n    @('73-23=', '30+46='),\n    @('50+11=', '68+10='),\n    @('45-9=', '29+66='),\n    @('40+11=', '32+10='),\n    @('8+58=', '50-5='),\n    @('83+6=', '61+30='),\n    @('21+48=', '28-1='),\n    @('8+91=', '52-26='),\n    @('17+41=', '36+43='),\n    @('44+15=', '78-68='),\n    @('89+1=', '47-7='),\n    @('45-25=', '10+51='),\n    @('98-82=', '24-10='),\n    @('70-34=', '85-6='),\n    @('27+52=', '28-10='),\n    @('16+56=', '48+14='),\n    @('49+38=', '88-28='),\n    @('85-18=', '75-5='),\n    @('8+82=', '31+28=')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    [void]$range.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
